$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-14, keeping header row 1 and data rows 2-5
$ws.Range("A6:E14").EntireRow.Delete()

# Update column widths: B -> 16, C -> 30 (offset by 5/6 to compensate for
# the engine's internal ColumnWidth <-> stored-width conversion)
$ws.Columns("B").ColumnWidth = 15.166666666666666
$ws.Columns("C").ColumnWidth = 29.166666666666668

# Phone/ID columns (D and E) must stay stored as text, not numbers, to
# match the original inline-string cell type. Temporarily force a text
# number format while assigning the values, then restore the default
# "Normal" style so no stray style attribute is left on the cells.
$dataRange = $ws.Range("D2:E5")
$dataRange.NumberFormat = "@"

# Row 2 data
$ws.Range("A2").Value = "Juan Carlos"
$ws.Range("B2").Value = "Pérez González"
$ws.Range("C2").Value = "jperez@sena.edu.co"
$ws.Range("D2").Value = "3001234567"
$ws.Range("E2").Value = "1234567890"

# Row 3 data
$ws.Range("A3").Value = "Jesus Andres"
$ws.Range("B3").Value = "Silva Plazas"
$ws.Range("C3").Value = "jsapp@sena.edu.co"
$ws.Range("D3").Value = "3182528515"
$ws.Range("E3").Value = "1055878001"

# Row 4 data
$ws.Range("A4").Value = "Gabriel"
$ws.Range("B4").Value = "Jesurum Rojas"
$ws.Range("C4").Value = "gabrieljesurumro@sena.edu.co"
$ws.Range("D4").Value = "3143887918"
$ws.Range("E4").Value = "3125435"

# Row 5 data
$ws.Range("A5").Value = "Hector"
$ws.Range("B5").Value = "Plaz Plaza"
$ws.Range("C5").Value = "hector@sena.edu.co"
$ws.Range("D5").Value = "3142884050"
$ws.Range("E5").Value = "107784658"

$dataRange.Style = "Normal"
